$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 8): Binary Search / Search Insert Position / notes
$ws.Range("A8").Value = "Binary Search"
$ws.Range("B8").Value = "Search Insert Position"
$ws.Range("C8").Value = "discrete binary search, lower_bound;"

# Match the highlight formatting used on the other "Name" column cells (copy from B6)
$ws.Range("B6").Copy()
$ws.Range("B8").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Leave the selection on the last edited cell, like the author's session
$ws.Range("C8").Select()
